# Auto-generated edit script: update market price columns (H-N) per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 332.27274
$ws.Range("I80").Value = 200.5
$ws.Range("J80").Value = 683.6667
$ws.Range("K80").Value = 601.5
$ws.Range("L80").Value = 2051.0001
$ws.Range("M80").Value = 396.5
$ws.Range("N80").Value = -4047.0001

$ws.Range("H83").Value = 332.27274
$ws.Range("I83").Value = 200.5
$ws.Range("J83").Value = 683.6667
$ws.Range("K83").Value = 1804.5
$ws.Range("L83").Value = 6153.0003
$ws.Range("M83").Value = 3187.5
$ws.Range("N83").Value = -16137.0003

$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496

$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480

$ws.Range("H116").Value = 4452.2104
$ws.Range("I116").Value = 2266.6667
$ws.Range("J116").Value = 5460.923
$ws.Range("K116").Value = 2266.6667
$ws.Range("L116").Value = 5460.923
$ws.Range("M116").Value = 1175.3333
$ws.Range("N116").Value = -12344.923

$ws.Range("H135").Value = 29239.584
$ws.Range("I135").Value = 32644.531
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 293800.779
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -291265.779
$ws.Range("N135").Value = -23070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25095.207
$ws.Range("I32").Value = 4715.094
$ws.Range("J32").Value = 133109.8
$ws.Range("K32").Value = 4715.094
$ws.Range("L32").Value = 133109.8
$ws.Range("M32").Value = -4428.094
$ws.Range("N32").Value = -133683.8

$ws.Range("H44").Value = 24839.6
$ws.Range("J44").Value = 24839.6
$ws.Range("L44").Value = 24839.6
$ws.Range("N44").Value = -25815.6

$ws.Range("H52").Value = 12800
$ws.Range("J52").Value = 12800
$ws.Range("L52").Value = 12800
$ws.Range("N52").Value = -13436

$ws.Range("H74").Value = 945.85297
$ws.Range("J74").Value = 808.125
$ws.Range("L74").Value = 808.125
$ws.Range("N74").Value = -2556.125

$ws.Range("H77").Value = 945.85297
$ws.Range("J77").Value = 808.125
$ws.Range("L77").Value = 4040.625
$ws.Range("N77").Value = -12776.625

$ws.Range("H86").Value = 40314
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 40314
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 40314
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -42686

$ws.Range("H89").Value = 40314
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 40314
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 120942
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -132798

$ws.Range("H105").Value = 27275
$ws.Range("J105").Value = 27275
$ws.Range("L105").Value = 27275
$ws.Range("N105").Value = -34263

$ws.Range("H119").Value = 22695
$ws.Range("J119").Value = 22695
$ws.Range("L119").Value = 22695
$ws.Range("N119").Value = -32371

$ws.Range("H132").Value = 1920.0625
$ws.Range("I132").Value = 1593.76
$ws.Range("J132").Value = 3085.4285
$ws.Range("K132").Value = 4781.28
$ws.Range("L132").Value = 9256.2855
$ws.Range("M132").Value = -2251.28
$ws.Range("N132").Value = -14316.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2347.95
$ws.Range("I134").Value = 2392.5789
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 7177.736699999999
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -4642.736699999999
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2749.36
$ws.Range("I31").Value = 2467.8823
$ws.Range("J31").Value = 2894.3635
$ws.Range("K31").Value = 2467.8823
$ws.Range("L31").Value = 2894.3635
$ws.Range("M31").Value = -2172.8823
$ws.Range("N31").Value = -3484.3635

$ws.Range("H34").Value = 2749.36
$ws.Range("I34").Value = 2467.8823
$ws.Range("J34").Value = 2894.3635
$ws.Range("K34").Value = 2467.8823
$ws.Range("L34").Value = 2894.3635
$ws.Range("M34").Value = -2265.8823
$ws.Range("N34").Value = -3298.3635

$ws.Range("H81").Value = 29900
$ws.Range("J81").Value = 29900
$ws.Range("L81").Value = 29900
$ws.Range("N81").Value = -31896

$ws.Range("H84").Value = 29900
$ws.Range("J84").Value = 29900
$ws.Range("L84").Value = 89700
$ws.Range("N84").Value = -99684

$ws.Range("H88").Value = 43995
$ws.Range("J88").Value = 43995
$ws.Range("L88").Value = 43995
$ws.Range("N88").Value = -44807

$ws.Range("H91").Value = 43995
$ws.Range("J91").Value = 43995
$ws.Range("L91").Value = 43995
$ws.Range("N91").Value = -46803

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 1751.4286
$ws.Range("J48").Value = 2129.0908
$ws.Range("L48").Value = 6387.2724
$ws.Range("N48").Value = -6887.2724

$ws.Range("H75").Value = 2359.2856
$ws.Range("J75").Value = 2669.1667
$ws.Range("L75").Value = 8007.500100000001
$ws.Range("N75").Value = -10003.5001

$ws.Range("H78").Value = 2359.2856
$ws.Range("J78").Value = 2669.1667
$ws.Range("L78").Value = 24022.5003
$ws.Range("N78").Value = -34006.5003

$ws.Range("H81").Value = 100001800
$ws.Range("I81").Value = 406.5
$ws.Range("J81").Value = 125002150
$ws.Range("K81").Value = 1219.5
$ws.Range("L81").Value = 375006450
$ws.Range("M81").Value = -96.5
$ws.Range("N81").Value = -375008696

$ws.Range("H84").Value = 100001800
$ws.Range("I84").Value = 406.5
$ws.Range("J84").Value = 125002150
$ws.Range("K84").Value = 3658.5
$ws.Range("L84").Value = 1125019350
$ws.Range("M84").Value = 1957.5
$ws.Range("N84").Value = -1125030582

$ws.Range("H114").Value = 2374.4614
$ws.Range("I114").Value = 705.4286
$ws.Range("J114").Value = 4321.6665
$ws.Range("K114").Value = 2116.2858
$ws.Range("L114").Value = 12964.9995
$ws.Range("M114").Value = 1137.7142
$ws.Range("N114").Value = -19472.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1448.6923
$ws.Range("I113").Value = 1508.1052
$ws.Range("J113").Value = 1287.4286
$ws.Range("K113").Value = 1508.1052
$ws.Range("L113").Value = 1287.4286
$ws.Range("M113").Value = 661.8948
$ws.Range("N113").Value = -5627.4286

$ws.Range("H126").Value = 5649.143
$ws.Range("I126").Value = 2641.5715
$ws.Range("J126").Value = 8656.714
$ws.Range("K126").Value = 7924.7145
$ws.Range("L126").Value = 25970.142
$ws.Range("M126").Value = -5454.7145
$ws.Range("N126").Value = -30910.142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3194.375
$ws.Range("I61").Value = 2040.6471
$ws.Range("J61").Value = 5996.2856
$ws.Range("K61").Value = 2040.6471
$ws.Range("L61").Value = 5996.2856
$ws.Range("M61").Value = -1838.6471
$ws.Range("N61").Value = -6400.2856

$ws.Range("H80").Value = 34777
$ws.Range("J80").Value = 34777
$ws.Range("L80").Value = 34777
$ws.Range("N80").Value = -37023

$ws.Range("H83").Value = 34777
$ws.Range("J83").Value = 34777
$ws.Range("L83").Value = 104331
$ws.Range("N83").Value = -115563

$ws.Range("H113").Value = 3194.375
$ws.Range("I113").Value = 2040.6471
$ws.Range("J113").Value = 5996.2856
$ws.Range("K113").Value = 2040.6471
$ws.Range("L113").Value = 5996.2856
$ws.Range("M113").Value = 129.3529000000001
$ws.Range("N113").Value = -10336.2856

$ws.Range("H116").Value = 21600
$ws.Range("J116").Value = 21600
$ws.Range("L116").Value = 21600
$ws.Range("N116").Value = -30778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H108").Value = 29900
$ws.Range("J108").Value = 29900
$ws.Range("L108").Value = 29900
$ws.Range("N108").Value = -37580

$ws.Range("H132").Value = 4580.6924
$ws.Range("I132").Value = 5194.1177
$ws.Range("K132").Value = 15582.3531
$ws.Range("M132").Value = -13052.3531

$ws.Range("H136").Value = 2271.0527
$ws.Range("I136").Value = 2013.0667
$ws.Range("J136").Value = 3238.5
$ws.Range("K136").Value = 6039.2001
$ws.Range("L136").Value = 9715.5
$ws.Range("M136").Value = -3489.2001
$ws.Range("N136").Value = -14815.5
